# Weekly update: a new daily price record for Haba (Vega Central Mapocho de
# Santiago) is inserted as row 406, pushing the previously-existing rows
# 406:431 down to 407:432 (dimension grows from A1:R431 to A1:R432).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 406, shifting rows 406:431 down to 407:432.
$ws.Rows.Item(406).Insert()

# Populate the new row 406 with the new record's data.
$ws.Cells.Item(406, 1).Value = 9
$ws.Cells.Item(406, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(406, 3).Value = "Metropolitana"
$ws.Cells.Item(406, 4).Value = 45267
$ws.Cells.Item(406, 5).Value = 13
$ws.Cells.Item(406, 6).Value = 100112026
$ws.Cells.Item(406, 7).Value = "Haba"
$ws.Cells.Item(406, 8).Value = "Sin especificar"
$ws.Cells.Item(406, 9).Value = "Primera"
$ws.Cells.Item(406, 10).Value = 70
$ws.Cells.Item(406, 11).Value = 11000
$ws.Cells.Item(406, 12).Value = 13000
$ws.Cells.Item(406, 13).Value = 12000
$ws.Cells.Item(406, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(406, 15).Value = "Región del Maule"
$ws.Cells.Item(406, 16).Value = 480
$ws.Cells.Item(406, 17).Value = 25
$ws.Cells.Item(406, 18).Value = "Hortaliza"
